# Fixed defect number 1: update Ticket Sales (Q) and Embarking (R) values
# for the station rows that previously held placeholder 0 values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 3;   Q = 35;  R = 9  }
    @{ Row = 10;  Q = 8;   R = 1  }
    @{ Row = 17;  Q = 89;  R = 40 }
    @{ Row = 23;  Q = 62;  R = 50 }
    @{ Row = 32;  Q = 76;  R = 75 }
    @{ Row = 40;  Q = 15;  R = 7  }
    @{ Row = 49;  Q = 100; R = 96 }
    @{ Row = 58;  Q = 56;  R = 49 }
    @{ Row = 66;  Q = 95;  R = 8  }
    @{ Row = 74;  Q = 18;  R = 4  }
    @{ Row = 78;  Q = 92;  R = 54 }
    @{ Row = 89;  Q = 19;  R = 15 }
    @{ Row = 97;  Q = 21;  R = 1  }
    @{ Row = 106; Q = 59;  R = 20 }
    @{ Row = 115; Q = 70;  R = 46 }
    @{ Row = 124; Q = 9;   R = 6  }
    @{ Row = 133; Q = 95;  R = 44 }
    @{ Row = 142; Q = 98;  R = 61 }
)

foreach ($u in $updates) {
    $ws.Range("Q$($u.Row)").Value = $u.Q
    $ws.Range("R$($u.Row)").Value = $u.R
}
